# Add a dark box (thicker outside border, thin inside gridlines) around the
# two "privileges" tables in the Driving Privileges template.
#
# wdLineStyleSingle = 1
# The runtime's Border.LineWidth is expressed in half of the OOXML `sz`
# units (eighths of a point), i.e. LineWidth = sz / 2.
#   sz=12 (1.5pt) outer frame  -> LineWidth = 6
#   sz=2  (0.25pt) inner grid  -> LineWidth = 1
#
# wdBorderTop = -1, wdBorderLeft = -2, wdBorderBottom = -3,
# wdBorderRight = -4, wdBorderHorizontal = -5, wdBorderVertical = -6

$d = $word.ActiveDocument

function Set-DarkBoxBorders($table) {
    foreach ($idx in -1, -2, -3, -4) {
        $border = $table.Borders.Item($idx)
        $border.LineStyle = 1
        $border.LineWidth = 6
    }
    foreach ($idx in -5, -6) {
        $border = $table.Borders.Item($idx)
        $border.LineStyle = 1
        $border.LineWidth = 1
    }
}

# Table 1: the case/defendant/suspension summary table.
Set-DarkBoxBorders($d.Tables.Item(1))

# Table 2: the employer/school driving-privileges table.
Set-DarkBoxBorders($d.Tables.Item(2))
